$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 468, shifting existing rows 468-568 down to 469-569
$ws.Rows.Item(468).Insert()

# Populate the newly inserted row 468 with the new record's data.
# Columns A,B,C,E,F,G,H,I,O,R are identical to the (now shifted) row below it,
# columns D,J,K,L,M,N,P,Q carry the new values from the edit.
$ws.Range("A468").Value2 = 7
$ws.Range("B468").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C468").Value2 = "Ñuble"
$ws.Range("D468").Value2 = 44889
$ws.Range("E468").Value2 = 16
$ws.Range("F468").Value2 = 100112020
$ws.Range("G468").Value2 = "Tomate"
$ws.Range("H468").Value2 = "Larga vida"
$ws.Range("I468").Value2 = "Primera"
$ws.Range("J468").Value2 = 400
$ws.Range("K468").Value2 = 20000
$ws.Range("L468").Value2 = 21000
$ws.Range("M468").Value2 = 20500
$ws.Range("N468").Value2 = '$/bandeja 18 kilos'
$ws.Range("O468").Value2 = "Región del Maule"
$ws.Range("P468").Value2 = 1139
$ws.Range("Q468").Value2 = 18
$ws.Range("R468").Value2 = "Hortaliza"
